# Team_Everyday_Attendence.xlsx — add 11 Aug 23 (row 10) and 12 Aug 23 (row 11)
# attendance rows, their review comments, and a new blank row 12 below them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 : 11-Aug-2023 (serial 45149) ---------------------------------
$ws.Range("A10").Value = 45149
$ws.Range("A10").NumberFormat = "d-mmm-yy"

$ws.Range("B10").Value = "ABSENT"
$ws.Range("C10").Value = "ABSENT"
$ws.Range("D10").Value = "PRESENT"
$ws.Range("E10").Value = "PRESENT"
$ws.Range("F10").Value = "PRESENT"
$ws.Range("G10").Value = "ABSENT"
$ws.Range("H10").Value = "ABSENT"
$ws.Range("I10").Value = "ABSENT"
$ws.Range("J10").Value = "ABSENT"
$ws.Range("K10").Value = "ABSENT"

# --- Row 11 : 12-Aug-2023 (serial 45150) ---------------------------------
$ws.Range("A11").Value = 45150
$ws.Range("A11").NumberFormat = "d-mmm-yy"

$ws.Range("B11").Value = "ABSENT"
$ws.Range("C11").Value = "PRESENT"
$ws.Range("D11").Value = "PRESENT"
$ws.Range("E11").Value = "PRESENT"
$ws.Range("F11").Value = "PRESENT"
$ws.Range("G11").Value = "ABSENT"
$ws.Range("H11").Value = "ABSENT"
$ws.Range("I11").Value = "ABSENT"
$ws.Range("J11").Value = "ABSENT"
$ws.Range("K11").Value = "ABSENT"

# --- Row 12 : new empty row started below, matching the date column style
$ws.Range("A12").NumberFormat = "d-mmm-yy"

# --- Reviewer comments on row 10 -----------------------------------------
[void]$ws.Range("B10").AddComment("LENOVO:`nMedical issue")
[void]$ws.Range("G10").AddComment("LENOVO:`nNo Responsse")
[void]$ws.Range("H10").AddComment("LENOVO:`nNo Responsse")
[void]$ws.Range("I10").AddComment("LENOVO:`nNo Responsse")
[void]$ws.Range("J10").AddComment("LENOVO:`nNo Responsse")
[void]$ws.Range("K10").AddComment("LENOVO:`nNo Responsse")

# --- Reviewer comments on row 11 -----------------------------------------
[void]$ws.Range("B11").AddComment("LENOVO:`nMedical issue")
[void]$ws.Range("G11").AddComment("LENOVO:`nNo Responsse")
[void]$ws.Range("H11").AddComment("LENOVO:`nNo Responsse`n")
[void]$ws.Range("I11").AddComment("LENOVO:`nNo Responsse")
[void]$ws.Range("J11").AddComment("LENOVO:`nNo Responsse")
[void]$ws.Range("K11").AddComment("LENOVO:`nNo Responsse")

# --- Selection moves to the new blank row, matching the saved UI state ---
[void]$ws.Range("A12").Select()
